$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting the existing rows 8-13 down to
# 9-14, to make room for a new "lucky spin wheel" panel event row.
$ws.Rows("8:8").Insert()

# Copy the row-7 formatting down into the freshly inserted row 8 so the
# per-cell styles (s="7"/"9") match what Excel itself would stamp on it.
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

# Populate the new row 8 with the "lucky spin wheel" panel event.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "|您可以通过幸运转盘面板，获取更多资源和道具。"
$ws.Range("E8").Value = "shop"

# The rows that were pushed down keep their B:E content, but the sequential
# Id column (A) needs renumbering to stay contiguous.
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9

# Grow the table to cover the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:E14"))

# Match the saved selection state from the edit.
$ws.Range("E8").Select()
